# Insert a new data row for "Jengibre" (Vega Central Mapocho de Santiago) as
# row 59, pushing the existing rows 59-103 down to 60-104 (weekly price
# update per the "Fruta / hortaliza, semanal" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59..end down by one to make room for the new weekly record.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with this week's record.
$ws.Range("A59").Value = 9
$ws.Range("B59").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = 44767
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = 100114007
$ws.Range("G59").Value = "Jengibre"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 700
$ws.Range("K59").Value = 17000
$ws.Range("L59").Value = 18000
$ws.Range("M59").Value = 17500
$ws.Range("N59").Value = "$/caja 13 kilos"
$ws.Range("O59").Value = "Perú"
$ws.Range("P59").Value = 1346
$ws.Range("Q59").Value = 13
$ws.Range("R59").Value = "Hortaliza"
